# Update the "取得日時" (retrieved-at) timestamps in column A
# from "2025-12-23 01:24:52" to "2025-12-23 01:57:31" for data rows 2-16
# on the "ランサーズ" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-12-23 01:24:52"
$newValue = "2025-12-23 01:57:31"

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
